# edit.ps1 -- apply ptb_eg5.3.docx revision (v1.1.0 docs update)
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: turn the single-space run right after
# "adds the percent of nonmissing observations." into a paragraph-
# internal line break (w:br type="textWrapping").
# ---------------------------------------------------------------------
$anchor = "adds the percent of nonmissing observations."
$full = $d.Content.Text
$pos = $full.IndexOf($anchor)
if ($pos -lt 0) { throw "anchor1 not found" }
$spaceStart = $pos + $anchor.Length
$spaceRange = $d.Range($spaceStart, $spaceStart + 1)
if ($spaceRange.Text -ne " ") { throw "unexpected text at spaceRange: [$($spaceRange.Text)]" }
$spaceRange.Text = ""
$spaceRange.InsertBreak(6)   # wdLineBreak = 6 -> <w:br w:type="textWrapping"/>

Write-Output "step1 done"

# ---------------------------------------------------------------------
# Step 2: rewrite the tail of the FirstParagraph paragraph:
#   "is specified as well to include a percentage sign."
# becomes a longer passage that also documents su_decimal()/miss_decimal().
# ---------------------------------------------------------------------
$oldTail = "is specified as well to include a percentage sign."
$newTail = "may be specified as well to include a percentage sign. When denominators or missing data summaries are included in the table the options su_decimal(#) and miss_decimal(#) can be used to independently control the number of decimal places reported for summary statistics and the percent of missing/nonmissing observations."

$full = $d.Content.Text
$pos = $full.IndexOf($oldTail)
if ($pos -lt 0) { throw "anchor2 not found" }
$tailRange = $d.Range($pos, $pos + $oldTail.Length)
if ($tailRange.Text -ne $oldTail) { throw "unexpected tailRange text: [$($tailRange.Text)]" }
$tailRange.Text = $newTail

Write-Output "step2 done"

# Re-apply the "Verbatim Char" style to the two inline option tokens that
# were just inserted as plain text.
foreach ($tok in @("su_decimal(#)", "miss_decimal(#)")) {
    $full = $d.Content.Text
    $p = $full.IndexOf($tok)
    if ($p -lt 0) { throw "token not found: $tok" }
    $tr = $d.Range($p, $p + $tok.Length)
    if ($tr.Text -ne $tok) { throw "mismatch for $tok : [$($tr.Text)]" }
    $tr.Style = "Verbatim Char"
}

Write-Output "step2b done"

# ---------------------------------------------------------------------
# Step 3: SourceCode example block updates.
# ---------------------------------------------------------------------

# 3a. "post `postname' (...)" header line: group "2" columns renamed "0".
$oldHeaderBit = '("N 2") ("Summary 2")'
$newHeaderBit = '("N 0") ("Summary 0")'
$full = $d.Content.Text
$p = $full.IndexOf($oldHeaderBit)
if ($p -lt 0) { throw "header bit not found" }
$hr = $d.Range($p, $p + $oldHeaderBit.Length)
if ($hr.Text -ne $oldHeaderBit) { throw "header mismatch: [$($hr.Text)]" }
$hr.Text = $newHeaderBit

Write-Output "step3a done"

# 3b. "pt_base age" command line gains trailing options.
$ageLine = '. pt_base age , post(`postname'') over(treat)  overall(last) over_grps(1, 0) type(cont) su_label(append) cat_col  n_analysis(cols cond %) order(group_over) per'
$full = $d.Content.Text
$p = $full.IndexOf($ageLine)
if ($p -lt 0) { throw "age line not found" }
$insPos = $p + $ageLine.Length
$insRange = $d.Range($insPos, $insPos)
$insRange.InsertBefore("  miss_decimal(2) su_decimal(0)")

Write-Output "step3b done"

# 3c. "pt_base qol" command line gains trailing options.
$qolLine = '. pt_base qol, post(`postname'') over(treat)  overall(last)  over_grps(1, 0) type(skew) su_label(append) cat_col  n_analysis(cols  cond %) order(group_over) per'
$full = $d.Content.Text
$p = $full.IndexOf($qolLine)
if ($p -lt 0) { throw "qol line not found" }
$insPos = $p + $qolLine.Length
$insRange = $d.Range($insPos, $insPos)
$insRange.InsertBefore("  miss_decimal(2) decimal(1)")

Write-Output "step3c done"

# ---------------------------------------------------------------------
# Sanity check: still exactly 3 paragraphs (Heading3 / FirstParagraph /
# SourceCode), nothing got split or merged unexpectedly.
# ---------------------------------------------------------------------
if ($d.Paragraphs.Count -ne 3) { throw "unexpected paragraph count: $($d.Paragraphs.Count)" }

Write-Output "all edits applied"
